$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: OA -> MMP62, DMSO -> MMAE, null (rows 9-11) -> DMSO
$ws.Range("A3").Value = "MMP62"
$ws.Range("A4").Value = "MMP62"
$ws.Range("A5").Value = "MMP62"

$ws.Range("A6").Value = "MMAE"
$ws.Range("A7").Value = "MMAE"
$ws.Range("A8").Value = "MMAE"

$ws.Range("A9").Value = "DMSO"
$ws.Range("A10").Value = "DMSO"
$ws.Range("A11").Value = "DMSO"

# Update selection to A11
$ws.Range("A11").Select()
